# Tracking Bot handout -- "Better extensions for tracking bot"
#
# Rewrites the three "Bonus" bullet paragraphs that used to read:
#   - Add an “evil sprite” that chases the drawing sprite and makes you lose if it touches it
#   - Add an explosion animation if the tracking bot touches brown
#   - Add the concept of lives
# into four bullets:
#   - Create a 3+ frame explosion animation when the tracking bot touches brown
#   - Add the concept of lives
#   - Alter the code so instead of following a line, the tracking bot “shies” away from a line
#   - Create a start screen that allows the user to select which mode (follow or shy) to use

$d = $word.ActiveDocument

$PKG_OPEN = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$PKG_CLOSE = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Find-ParaIndex($doc, $needle) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text.Contains($needle)) {
            return $idx
        }
    }
    return -1
}

function Replace-ParaBody($doc, $paraIndex, $innerXml) {
    # Replace everything in the paragraph *except* the trailing paragraph
    # mark with freshly authored OOXML; $innerXml must carry its own
    # <w:pPr> (we pass the original one back in) since InsertXML on a
    # collapsed range swaps in a brand new <w:p>.
    $p = $doc.Paragraphs($paraIndex)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = ""
    $r2 = $doc.Paragraphs($paraIndex).Range
    $r2.MoveEnd(1, -1) | Out-Null
    $xml = $PKG_OPEN + "<w:body>" + $innerXml + "</w:body>" + $PKG_CLOSE
    $r2.InsertXML($xml)
}

function Insert-ParaAfter($doc, $paraIndex, $innerXml) {
    # Insert a brand new paragraph right after $paraIndex.
    $p = $doc.Paragraphs($paraIndex)
    $r = $p.Range
    $r.InsertParagraphAfter()
    $newPara = $doc.Paragraphs($paraIndex + 1)
    $r2 = $newPara.Range
    $r2.MoveEnd(1, -1) | Out-Null
    $xml = $PKG_OPEN + "<w:body>" + $innerXml + "</w:body>" + $PKG_CLOSE
    $r2.InsertXML($xml)
}

$GFONT = '<w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>'

# Locate all four target paragraphs *before* any edits are made, so later
# text changes can never shadow an earlier Find-ParaIndex lookup.
$iEvil = Find-ParaIndex $d 'evil sprite'
$iExplosion = Find-ParaIndex $d 'explosion animation if'
$iLives = Find-ParaIndex $d 'Add the concept of lives'

# ---------------------------------------------------------------------
# 1) "evil sprite" paragraph -> "Create a 3+ frame explosion animation
#    when the tracking bot touches brown"
# ---------------------------------------------------------------------
$p1Inner = ('<w:p><w:pPr><w:ind w:firstLine="540"/>' + $GFONT + '</w:pPr>' +
  '<w:r><w:t xml:space="preserve">□ </w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t>Create</w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $GFONT + '<w:t>a</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve"> 3+ frame</w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve"> explosion animation </w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t>when</w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve"> the tracking bot touches brown</w:t></w:r>' +
  '</w:p>')
Replace-ParaBody $d $iEvil $p1Inner

# ---------------------------------------------------------------------
# 2) "explosion animation" paragraph -> "Add the concept of lives"
# ---------------------------------------------------------------------
$p2Inner = ('<w:p><w:pPr><w:ind w:firstLine="540"/>' + $GFONT + '</w:pPr>' +
  '<w:r><w:t xml:space="preserve">□ </w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve">Add </w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve">the concept of </w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t>lives</w:t></w:r>' +
  '</w:p>')
Replace-ParaBody $d $iExplosion $p2Inner

# ---------------------------------------------------------------------
# 3) "lives" paragraph -> "Alter the code so instead of following a
#    line, the tracking bot “shies” away from a line"
# ---------------------------------------------------------------------
$p3Inner = ('<w:p><w:pPr><w:ind w:firstLine="540"/>' + $GFONT + '</w:pPr>' +
  '<w:r><w:t>□</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t>Alter the code so instead of following a line, the tracking bot “shies” away from a line</w:t></w:r>' +
  '</w:p>')
Replace-ParaBody $d $iLives $p3Inner

# ---------------------------------------------------------------------
# 4) brand new paragraph after it -> "Create a start screen that
#    allows the user to select which mode (follow or shy) to use"
# ---------------------------------------------------------------------
$p4Inner = ('<w:p><w:pPr><w:ind w:firstLine="540"/>' + $GFONT + '</w:pPr>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve">        </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t xml:space="preserve">□ </w:t></w:r>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve"> Create</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r>' + $GFONT + '<w:t xml:space="preserve"> a start screen that allows the user to select which mode (follow or shy) to use</w:t></w:r>' +
  '</w:p>')
Insert-ParaAfter $d $iLives $p4Inner

# ---------------------------------------------------------------------
# 5) Section bottom margin 1440 -> 1341 twips (67.05pt)
# ---------------------------------------------------------------------
$d.Sections(1).PageSetup.BottomMargin = 1341 / 20.0
